$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "Penalty Issued" -------------------------------------
# Header cell I4 should look exactly like the other header cells on row 4
# (bold, centered, wrapped) - copy the format from H4, then set its text.
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").Value = "Penalty Issued"

# Data placeholder cell I6 (plain, default formatting), mirrors the other
# merge-field placeholders used throughout the template.
$ws.Range("I6").Value = "{d.Reg[i].PenaltyIssued}"

# Spacer cells I3 / I5 get a thin-bottom-border-only style (no fill, no
# font change, no alignment) - build it from scratch via the Borders API.
foreach ($addr in @("I3", "I5")) {
    $border = $ws.Range($addr).Borders.Item(9)
    $border.Color = -16777216
    $border.LineStyle = 1
}

# Column I width to match the rest of the report's header row.
$ws.Columns.Item(9).ColumnWidth = 10.45

# Restore the cursor to where the author left it after adding the column.
$ws.Range("G16").Select() | Out-Null
